# Auto-applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. '63.845.47', '0.160').
# Force Text number-format before assigning so Excel doesn't coerce them to
# doubles (which would strip significant trailing/leading zeros), then drop
# back to the default 'Normal' style so no stray formatting is left behind.
function Set-TextValue([string]$cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.845.47"
$ws.Range("E2").Value = "  -0.46%  "
Set-TextValue "D3" "2.751.81"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "572.48"
$ws.Range("E5").Value = "  -1.65%  "
Set-TextValue "D6" "156.96"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("E9").Value = "  -3.96%  "
Set-TextValue "D10" "0.160"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D11" "5.63"
$ws.Range("E11").Value = "  -16.85%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D12" "0.381"
$ws.Range("E12").Value = "  -3.09%  "
Set-TextValue "D13" "3.238.42"
$ws.Range("E13").Value = "  +0.23%  "
Set-TextValue "D14" "26.38"
$ws.Range("E14").Value = "  -3.20%  "
Set-TextValue "D15" "63.502.08"
$ws.Range("E15").Value = "  -0.87%  "
Set-TextValue "D16" "0.0000150"
$ws.Range("E16").Value = "  -3.10%  "
Set-TextValue "D17" "2.757.95"
$ws.Range("E17").Value = "  -0.43%  "
Set-TextValue "D18" "12.09"
$ws.Range("E18").Value = "  -0.20%  "
Set-TextValue "D19" "4.79"
$ws.Range("E19").Value = "  -3.08%  "
Set-TextValue "D20" "354.03"
$ws.Range("E20").Value = "  -2.69%  "
Set-TextValue "D21" "6.73"
$ws.Range("E21").Value = "  -4.92%  "
$ws.Range("E22").Value = "  +0.01%  "
Set-TextValue "D23" "0.533"
$ws.Range("E23").Value = "  -0.36%  "
Set-TextValue "D24" "65.08"
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("E26").Value = "  +0.12%  "
Set-TextValue "D27" "8.39"
$ws.Range("E27").Value = "  -2.89%  "
Set-TextValue "D28" "0.0₃0900"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("E29").Value = "  -4.77%  "
Set-TextValue "D30" "6.99"
$ws.Range("E30").Value = "  -3.10%  "
Set-TextValue "D31" "169.11"
$ws.Range("E31").Value = "  -3.38%  "
Set-TextValue "D32" "1.20"
$ws.Range("E32").Value = "  -5.91%  "
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("E34").Value = "  +0.14%  "
Set-TextValue "D35" "4.85"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  -2.59%  "
Set-TextValue "D38" "0.976"
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("E39").Value = "  +4.31%  "
Set-TextValue "D40" "4.13"
$ws.Range("E40").Value = "  -3.81%  "
Set-TextValue "D41" "327.02"
$ws.Range("E41").Value = "  -4.67%  "
Set-TextValue "D42" "39.01"
$ws.Range("E42").Value = "  -1.13%  "
Set-TextValue "D43" "21.41"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("E44").Value = "  -2.89%  "
Set-TextValue "D45" "21.27"
$ws.Range("E45").Value = "  -4.21%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D46" "0.0253"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "134.85"
$ws.Range("E47").Value = "  -2.90%  "
Set-TextValue "D48" "0.624"
$ws.Range("E48").Value = "  -4.30%  "
Set-TextValue "D49" "0.101"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("E50").Value = "  +0.24%  "
Set-TextValue "D51" "11.03"
$ws.Range("E51").Value = "  +0.15%  "
